# Apply "Feb 19th" offsets update to ImportStaff.xlsx
# - Staff Vitals sheet: drop the old non-prefixed "Last Name"/"Salary" helper
#   columns (keep "Position"), shift everything over, and resort/relabel the
#   "Staff Vitals - *" columns (adds LASTNAME).
# - Staff Attributes sheet: drop the old non-prefixed helper columns (keep
#   "Current Team"), and resort the "Staff Attributes - *" columns.
# - Staff Style sheet: rename/insert a few proficiency + style columns.

$wb = $excel.ActiveWorkbook

function Set-RowValues {
    # NOTE: use positional parameters here - passing arrays through named
    # parameters (e.g. "-Values $arr") loses the array contents in this
    # runtime, so callers below invoke this positionally: Set-RowValues $sheet $values
    param($Sheet, $Values)

    # Clear out the full previously-used range on row 1, then write the
    # new header row so stale trailing cells don't linger.
    $usedCols = $Sheet.UsedRange.Columns.Count
    if ($usedCols -gt 0) {
        $Sheet.Range($Sheet.Cells.Item(1, 1), $Sheet.Cells.Item(1, $usedCols)).Clear() | Out-Null
    }

    for ($i = 0; $i -lt $Values.Count; $i++) {
        $Sheet.Cells.Item(1, $i + 1).Value = $Values[$i]
    }
}

# ---------------------------------------------------------------------
# Sheet: Staff Vitals
# ---------------------------------------------------------------------
$wsVitals = $wb.Worksheets.Item("Staff Vitals")

$vitalsValues = @(
    "Face ID",
    "Height",
    "Position",
    "Staff Vitals - ARM_SCALE",
    "Staff Vitals - BODYLENGTH",
    "Staff Vitals - BODY_SHAPE",
    "Staff Vitals - CURRENT_TEAM",
    "Staff Vitals - EYE_COLOR",
    "Staff Vitals - FIRSTNAME",
    "Staff Vitals - GENDER",
    "Staff Vitals - HAIR_LENGTH",
    "Staff Vitals - HAND_SCALE",
    "Staff Vitals - HEIGHT_CM",
    "Staff Vitals - LASTNAME",
    "Staff Vitals - LOWER_SCALE",
    "Staff Vitals - NECK_HEAD_SCALE",
    "Staff Vitals - PERSONALITY",
    "Staff Vitals - POSITION",
    "Staff Vitals - SALARY",
    "Staff Vitals - SHOULDERWIDTH",
    "Staff Vitals - SKINCOLOR",
    "Staff Vitals - SKINTYPE",
    "Staff Vitals - UNIQUE_PHOTO_ID",
    "Staff Vitals - WINGSPAN_CM",
    "Staff Vitals - YEARS_IN_LEAGUE",
    "Staff Vitals - YEARS_LEFT"
)

Set-RowValues $wsVitals $vitalsValues

# ---------------------------------------------------------------------
# Sheet: Staff Attributes
# ---------------------------------------------------------------------
$wsAttributes = $wb.Worksheets.Item("Staff Attributes")

$attributesValues = @(
    "Current Team",
    "Staff Attributes - BUSINESS",
    "Staff Attributes - CONTRACTS",
    "Staff Attributes - DEFENSE",
    "Staff Attributes - MAX_BUSINESS",
    "Staff Attributes - MAX_CONTRACTS",
    "Staff Attributes - MAX_DEFENSE",
    "Staff Attributes - MAX_OFFENSE",
    "Staff Attributes - MAX_SCOUTING",
    "Staff Attributes - MAX_TRADING",
    "Staff Attributes - MAX_TRAINING",
    "Staff Attributes - OFFENSE",
    "Staff Attributes - POTENTIAL",
    "Staff Attributes - SCOUTING",
    "Staff Attributes - TRADING",
    "Staff Attributes - TRAINING"
)

Set-RowValues $wsAttributes $attributesValues

# ---------------------------------------------------------------------
# Sheet: Staff Style
# ---------------------------------------------------------------------
$wsStyle = $wb.Worksheets.Item("Staff Style")

$styleValues = @(
    "Grit & Grind Proficiency",
    "Pace &Space Proficiency",
    "Perimeter Centric Proficiency",
    "Post Centric Proficiency",
    "Staff Style - ACTIVE_SYSTEM",
    "Staff Style - BALANCED_PROFICIENCY",
    "Staff Style - DEFENSE_PROFICIENCY",
    "Staff Style - GUARDS_VS_FORWARDS",
    "Staff Style - INSIDE_VS_OUTSIDE",
    "Staff Style - OFFENSE_VS_DEFENSE",
    "Staff Style - PERIMETER_CENTERIC_PROFICIENCY",
    "Staff Style - PREFERRED_SYSTEM",
    "Staff Style - SEVEN_SECONDS_PROFICIENCY",
    "Staff Style - STYLE_N#1",
    "Staff Style - STYLE_N#2",
    "Staff Style - STYLE_N#3",
    "Triangle Proficiency"
)

Set-RowValues $wsStyle $styleValues
